$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 16 per-play logs ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "2 -1 19 2 13 2 0 1 1 1 4 4 5 3 21 10 6 1 2 -2 14 14 2 3 2 1 7 4 4 7 4 1 5 0 4 1 25 5 8 1 3 4 4 0 2 1 6 0 27 2 2 0 2 1 3 -2 8 1 3 5 8 0 8 0 0 9 5 2 0 -2 2 3 5 1 0 -3 12 1 7 2 1 7 4 5 -1 0 10 7 1 7 2 3 8 -2 3 5 3 5 3 17 3 3 1 54 1 3 0 1 -2 -1 3 5 9 2 -1 1 5 2 9 5 2 1 6 19 5 5 0 0 3 3 3 8 7 -1 1 8 13 3 3 0 7 0 -2 2 7 3 2 3 3 2 -5 8 1 7 -5 5 -3 0 2 -2 1 11 -1 3 10 8 4 4 2 8 16 5 5 4 9 1 2 2 -4 2 5 2 2 17 6 -1 2 4 4 1 5 14 5 11 16 1 9 2 13 4 0 0 3 0 2 2 2 3 -1 17 9 6 1 4 3 2 3 0 4 1 6 4 0 2 8 -1 1 4 0 1 0 1 1 0 8 3 7 4 2 4 1 0 -2 3 1 4 19 11 4 4 2 8 3 3 5 3 8 4 2 0 4 15 4 2 1 3 5 2 5 3 3 2 3 0 1 3 2 5 1 3 5 0 3 4 11 2 1 3 3 4 3 6 2 7 3 7 8 3 -1 3 4 5 2 12 4 0 1 2 7 11 4 7 2 0 10 -1 4 11 6 7 14 4 9 4 -4 3 3 0 0 1 16 0 11 3 18 4 3 4 5 3 38 4 1 6 4 6 -1 5 4 6 3 2 13 8 5 2 10 2 4 9 20 7 7 1 1 3 16 15 -1 -4 3 -4 1 3 6 1 10 1 1 4 2 8 5 2 6 3 3 8 3 1 12 2 17 3 2 0 26 11 1 4 14 2 6 4 14 4 0 9 -2 2 1 2 3 12 4 0 4 0 1 1 -1 0 1 -4 5 5 8 4 5 3 -1 1 3 11 2 -2 2 4 9 2 -1 2 4 6 6 3 1 3 -4 2 3 12 1 1 11 4 9 5 4 3 1 2 4 13 4 9 0 3 3 0 4 3 7 3 1 0 3 8 1 6 2 1 1 -5 1 0 9 5 -2 1 6 1 1 8 3 7 6 3 1 5 12 3 14 1 -3 3 4 0 1 3 2 2 11 14 4 5 0 4 1 1 1 -2 5 4 2 2 -2 5 0 2 7 8 1 7 9 3 0 0 4 28 15 4 2 5 21 -1 6 16 1 12 14 42 4 1 4 1 0 -5 0 4 -4 2 6 5 3 0 1 -2 8 3 3 -1 1 21 9 -2 3 -1 5 0 12 4 19 2 57 5 13 7 2 1 0 5 4 -4 0 0 2 4 0 3 7 1 8 4 2 4 2 6 11 19 3 2 10 1 2 8 4 5 4 6 8 2 4 3 0 2 -2 2 11 1 0 5 9 3 5 6 0 5 1 13 35 5 9 3 4 4 2 3 5 4 -3 7 3 3 6 -3 8 10 0 7 6 -3 5 5 0 27 6 8 0 0 -3 1 0 1 4 5 3 0 6 0 2 8 1 16 5 0 2 1 2 -2 3 2 7 1 3 2 3 14 5 -2 4 6 4 4 5 10 2 3 1 12 2 2 7 2 3 4 5 5 0 -3 4 4 7"
$wsYDS.Range("C2").Value = "6 3 0 13 3 6 5 -1 16 3 0 8 8 1 -2 8 5 1 7 8 10 2 3 0 6 20 2 5 12 4 -2 7 4 5 5 2 8 11 4 5 13 75 18 4 11 3 5 3 9 4 4 4 7 1 3 14 3 3 -1 0 12 5 4 0 5 7 1 8 1 1 0 11 3 4 1 5 1 -6 15 4 5 13 -1 2 2 14 1 3 6 2 7 5 4 8 3 0 3 7 5 9 0 1 3 6 6 12 10 6 1 1 5 2 9 5 6 1 1 8 4 4 11 -1 -2 -2 2 1 6 5 0 -1 2 6 1 1 3 5 5 -2 0 2 6 8 9 1 4 1 2 2 0 0 0 1 2 1 2 3 1 2 5 13 2 9 -6 -3 3 8 1 2 10 2 1 -2 4 2 5 -1 3 0 -2 3 3 16 6 3 2 3 4 4 1 1 0 -4 0 1 4 8 6 12 5 5 2 3 2 3 14 13 5 5 9 11 5 0 7 1 5 13 14 4 5 12 4 2 0 1 4 8 70 5 9 4 4 13 3 14 18 1 1 5 4 0 1 5 4 3 7 10 3 6 27 2 2 -3 1 2 3 1 12 3 2 4 2 2 5 1 5 5 3 2 6 1 2 1 4 2 1 9 4 3 0 5 4 7 3 0 4 21 9 5 3 4 3 0 6 2 3 2 3 0 12 8 -1 4 11 5 4 2 5 4 5 5 1 4 7 1 1 4 8 -1 3 12 13 8 13 4 1 2 4 10 5 7 3 0 6 4 1 13 2 2 2 -1 3 9 -3 3 3 3 5 3 6 2 5 1 2 6 4 1 4 11 0 4 6 8 -1 9 3 8 9 7 2 4 4 1 4 1 2 -4 4 5 9 5 3 12 8 13 6 7 17 2 8 10 5 1 -2 1 5 4 5 33 3 3 4 -1 3 3 4 8 9 4 -1 1 2 8 4 2 3 4 6 4 14 12 9 2 5 3 1 9 1 2 -1 1 4 4 -4 5 2 1 6 4 1 5 5 20 1 2 5 2 4 1 4 7 4 19 5 4 1 0 0 1 5 -1 0 11 9 1 3 3 38 6 11 3 3 3 2 10 2 7 3 2 1 7 -1 0 0 6 0 1 1 3 6 7 1 3 9 4 1 4 8 6 5 4 3 8 3 1 -4 5 0 3 4 1 8 -1 3 5 2 -2 2 5 0 31 6 2 0 4 8 8 1 9 0 5 7 8 5 12 3 9 4 1 4 4 10 1 4 -1 5 1 11 9 1 4 10 2 3 2 10 1 4 3 3 24 4 -1 5 7 2 23 8 -3 4 -1 4 4 -2 -2 9 4 0 5 5 11 3 6 0 3 3 3 0 5 48 2 4 -1 3 2 2 1 2 1 0 1 0 9 8 5 2 3 0 2 4 6 6 14 4 3 5 1 7 6 3 4 6 2 4 5 18 4 2 1 -2 6 4 4 3 5 10 0 3 1 1 1 4 0 -2 3 8 4 4 0 5 7 11 3 10 5 1 1 8 6 10 4 3 19 1 12 3 14 4 -2 20 4 7 3 10 4 3 21 4 3 5 2 2 5 3 3 11 2 5 5 0 2 9 3 12 26 1 7 4 11 1 0 6 8 4 2 3 3 1 3 2 4 0 -1 14 2 1 1 12 4 16 1 1 5 13 2 0 6 2 5 10 0 0 9 7 20 5 3 11 3 3 4 11 15 5 3 3 7 2 0 0 0 1 2 4 8 4 0 2 6 0 0 3 3 5 4 5 12 -1 9 1 0 1 7 5 1 4 9 3 5 -4 3 10 -1 2 5 1 8 11 4 8 0 1 3 4 4 2 5 1 11 5 3 1 13 9 1 5 6 3 3 9 4 1 3 2 2 1 3 3 -1 6 16 5 5 5 9 14 -1 1 1 6 8 3 4 -1 10 5 1 0 6 2 5 1 22 0 3 23 6 1 2 1 6 3 0 4 12 3 2 2 3 2 1 0 8 3 2 9 3 3 -2 3 0"
$wsYDS.Range("B3").Value = "20 26 18 21 7 7 16 24 15 14 16 4 9 4 9 4 5 11 5 3 32 5 29 6 9 8 14 4 13 0 10 12 19 4 9 8 9 17 8 24 36 18 11 15 0 1 20 11 10 19 15 1 19 27 14 8 7 8 10 20 11 20 50 7 15 14 31 16 0 7 9 8 1 4 17 1 9 3 16 11 18 36 1 2 3 6 19 48 5 6 16 11 6 12 3 3 16 23 1 13 36 1 9 10 18 20 16 8 29 14 4 18 15 6 7 -2 13 22 29 11 25 6 11 10 16 -2 6 15 15 5 14 23 9 17 10 10 17 2 73 4 13 24 7 6 1 5 -1 13 0 9 6 4 10 5 12 5 20 8 11 15 9 16 14 9 12 9 19 28 2 12 11 8 13 9 11 5 55 4 9 27 19 8 6 2 4 8 26 4 5 18 15 5 0 13 9 5 10 10 9 7 6 21 9 11 5 1 4 19 5 5 35 8 14 18 3 3 10 6 51 17 16 14 12 12 20 2 6 3 7 9 10 4 2 10 9 4 11 16 14 7 7 6 5 21 19 5 4 15 14 14 3 25 49 5 22 15 12 7 9 24 9 23 6 5 14 22 17 5 13 25 21 5 8 16 15 1 4 16 17 8 23 10 10 10 15 8 16 -1 5 6 12 13 2 4 1 23 5 0 6 4 13 11 12 2 6 5 8 2 5 44 7 6 39 12 14 7 36 2 11 7 21 11 5 11 15 2 7 2 2 7 4 19 11 2 6 35 22 15 3 10 6 5 20 2 9 7 -1 9 23 43 3 15 20 38 2 7 26 6 31 18 21 1 -1 16 6 10 22 6 -8 12 12 1 3 6 2 6 7 10 -1 7 11 6 43 3 23 5 5 3 10 5 6 5 3 8 7 9 11 13 3 20 2 20 25 8 46 9 5 6 4 5 11 4 10 8 20 1 6 15 7 2 5 3 5 7 7 5 24 16 7 11 3 4 5 8 1 15 10 5 12 11 14 19 9 2 -3 11 9 6 24 19 22 7 21 5 6 33 17 2 20 11 22 5 9 4 11 9 24 13 25 17 3 8 9 9 10 8 13 4 13 8 4 4 6 21 8 14 5 2 17 2 8 7 18 21 10 0 7 10 4 6 5 3 5 -2 5 4 33 2 9 5 3 15 7 4 3 7 10 3 7 8 10 14 15 4 63 -1 21 12 15 8 20 14 9 7 5 37 6 6 5 4 3 6 5 17 2 18 7 5 8 -1 5 7 7 2 10 5 9 18 4 34 0 8 6 11 11 7 9 17 5 2 -3 7 5 19 0 5 5 9 4 1 30 8 12 9 -1 5 3 2 9 -1 10 0 4 -1 6 24 -4 0 3 6 16 39 9 10 5 4 -6 9 2 0 -1 17 8 0 11 17 15 7 9 9 25 9 28 23 6 5 27 7 6 24 7 8 14 9 5 6 13 9 9 9 10 8 1 13 7 7 2 19 7 5 20 11 -2 9 9 11 13 4 6 14 8 14 7 2 4 5 15 5 13 5 13 5 37 5 15 14 8 22 16 6 5 8 11 -2 13 0 4 8 -3 12 2 14 20 8 6 8 2 6 4 24 2 23 5 4 8 5 8"
$wsYDS.Range("C3").Value = "13 12 10 19 8 11 12 22 19 3 0 18 22 3 5 2 7 18 13 16 1 27 15 4 22 11 4 10 7 17 8 8 11 30 14 25 41 9 -6 0 30 6 26 4 15 -3 13 18 11 11 15 10 9 21 14 6 15 12 11 10 4 12 16 18 19 8 19 25 10 2 20 15 7 4 12 -2 9 3 19 5 35 3 9 2 3 2 10 3 7 4 16 0 8 3 51 2 0 9 5 24 8 6 14 14 12 -2 7 16 7 9 19 11 27 11 17 4 13 13 4 9 5 17 8 8 5 28 7 12 6 9 8 16 16 24 5 5 8 4 -2 8 6 15 21 22 14 16 7 4 8 29 18 3 13 3 20 11 22 22 30 9 17 8 10 29 22 8 35 1 7 22 19 9 7 3 7 0 8 6 4 10 6 8 7 3 13 6 11 1 20 19 14 4 11 -1 13 8 8 23 27 5 26 5 11 11 11 9 12 8 14 52 10 2 4 6 3 8 7 15 7 18 20 17 -3 11 12 10 -2 7 11 10 11 14 2 33 12 37 27 -4 10 38 8 19 18 9 40 34 16 17 21 10 13 10 7 10 15 19 9 13 -1 6 4 21 6 4 11 7 9 16 1 14 5 4 4 15 56 10 3 5 12 14 11 4 8 7 3 14 14 -1 7 21 6 5 15 11 4 7 29 6 6 6 19 10 10 75 7 14 14 16 5 5 5 -1 20 14 7 24 2 12 3 8 33 33 47 -2 27 7 16 10 12 12 13 7 20 22 7 7 15 11 16 12 25 5 30 7 4 5 10 14 6 13 15 9 25 11 19 14 4 20 9 8 4 -3 14 10 40 15 38 15 7 14 30 13 8 5 23 7 5 29 8 7 23 7 14 18 35 8 79 16 2 15 12 11 18 19 4 5 6 9 13 1 10 50 22 14 3 9 5 22 11 11 0 15 8 0 5 16 11 6 21 37 19 29 3 9 41 7 24 5 18 4 4 21 64 8 8 28 32 6 27 7 37 14 4 17 6 3 6 14 22 2 15 4 9 2 11 6 3 14 9 7 16 8 21 6 19 14 24 5 5 1 9 4 3 7 9 34 8 7 23 11 53 14 2 7 16 29 4 5 1 5 5 2 11 11 27 3 22 10 4 9 2 7 13 22 5 14 59 14 5 6 18 5 7 16 0 15 13 17 19 43 15 6 11 4 9 9 -3 7 17 -2 1 7 20 -4 14 6 3 6 5 14 10 9 6 39 10 6 8 9 7 10 13 9 21 20 3 7 17 6 20 3 5 7 9 37 -1 3 11 13 8 7 12 33 2 17 52 17 16 8 9 0 10 19 11 16 12 7 13 12 7 19 0 18 -2 -2 34 6 7 14 6 2 16 15 11 5 19 10 48 7 16 4 21 14 16 5 14 5 -2 11 3 13 32 8 9 8 6 7 7 4 8 7 9 10 16 5 5 4 6 -1 14 13 2 5 24 10 15 0 3 9 20 19 27 16 5 7 6 7 13 11 26 -1 13 9 15 10 8 17 15 9 18 35 5 7 13 7 19 12 4"

# --- OFF sheet: cumulative totals updated for Week 16 ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("B2").Value = 7
$wsOFF.Range("C2").Value = 344
$wsOFF.Range("E2").Value = 23
$wsOFF.Range("F2").Value = 120
$wsOFF.Range("G2").Value = 94
$wsOFF.Range("I2").Value = 13
$wsOFF.Range("J2").Value = 47
$wsOFF.Range("O2").Value = 55
$wsOFF.Range("P2").Value = 27
$wsOFF.Range("B3").Value = 19
$wsOFF.Range("C3").Value = 355
$wsOFF.Range("E3").Value = 48
$wsOFF.Range("F3").Value = 186
$wsOFF.Range("G3").Value = 85
$wsOFF.Range("H3").Value = 52
$wsOFF.Range("I3").Value = 117
$wsOFF.Range("J3").Value = 98
$wsOFF.Range("L3").Value = 550
$wsOFF.Range("M3").Value = 360
$wsOFF.Range("Q3").Value = 1006

# --- DEF sheet: cumulative totals updated for Week 16 ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 448
$wsDEF.Range("F2").Value = 130
$wsDEF.Range("G2").Value = 118
$wsDEF.Range("J2").Value = 62
$wsDEF.Range("N2").Value = 22
$wsDEF.Range("O2").Value = 28
$wsDEF.Range("P2").Value = 17
$wsDEF.Range("C3").Value = 299
$wsDEF.Range("E3").Value = 51
$wsDEF.Range("F3").Value = 198
$wsDEF.Range("G3").Value = 79
$wsDEF.Range("H3").Value = 56
$wsDEF.Range("I3").Value = 95
$wsDEF.Range("J3").Value = 92
$wsDEF.Range("L3").Value = 513
$wsDEF.Range("M3").Value = 346
$wsDEF.Range("Q3").Value = 1001

# --- ST sheet: cumulative totals + per-game logs updated for Week 16 ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 132
$wsST.Range("D2").Value = 117
$wsST.Range("F2").Value = 127
$wsST.Range("G2").Value = 119
$wsST.Range("J2").Value = 59
$wsST.Range("K2").Value = 55
$wsST.Range("B3").Value = 52
$wsST.Range("B4").Value = "71 73 69 55 58 69 43 62 64 58 60 61 64 46 68 61 67 66 66 64 60 68 67 65 66 65 57 61 58 61 57 64 65 67 55 47 64 45 42 59 62 65 65 60 73 73 70 64 61 67 62 59 59 62 64 62 61 67 60 64 60 64 64 66 63 64 66 62 60 66 70 67 61 64 61 65 65 69 66 68"
$wsST.Range("B5").Value = "29 36 45 15 19 32 9 18 24 16 14 13 16 7 20 16 20 19 16 23 19 24 31 27 22 13 18 26 19 46 24 17 21 24 0 10 45 12 5 5 20 28 17 13 30 30 21 29 17 33 28 16 19 16 19 41 33 33 14 24 13 17 17 17 17 24 23 14 21 26 25 19 15 21 44 17 31 24 23 32"
$wsST.Range("B6").Value = "31 35 27 26 22 25 26 23 17 29 33 11 11 2 23 23 16 28 16 10 23 13 9 24 71 32 28 13 43 27 15 70 24 27 19 0 33 21 24 23 21 26 20 23 28 18 27 26 30 41 29 18 29 38 28 28 13 22 23 35 27 22 30 17 25"
$wsST.Range("D3").Value = "53 55 35 45 57 50 52 67 65 54 60 43 52 56 57 53 50 40 67 51 59 12 59 41 49 58 37 56 42 51 51 47 56 38 37 33 43 46 49 39 50 45 45 42 43 39 52 36 56 51 42 52 46 52 51 32 58 43 57 38 57 52 55 48 25 38 57 59 46 56 54 56 43 67 47 47 44 49 64 63 64 49 36 39 43 58 52 61 52 39 42 49 46 43 33 47 45 43 61 32 59 70 41 45 43 52 49 62 47 57 54 52 40 55 38 52 57"
$wsST.Range("D4").Value = "0 0 0 0 0 0 0 6 0 0 0 0 11 22 0 0 0 0 0 0 0 55 12 0 0 0 0 -3 0 0 0 0 8 0 0 0 7 14 0 0 0 0 5 0 0 0 2 0 0 6 5 -2 8 0 0 3 13 7 7 0 11 2 0 0 0 0 13 29 0 0 6 0 0 0 2 0 0 0 15 11 0 0 0 0 0 12 13 13 0 0 5 13 0 10 0 1 0 12 0 0 0 0 22 5 13 0 -1 16 2 5 2 11 6 12 0 13 0"
$wsST.Range("D5").Value = "0 0 0 16 0 0 0 0 0 19 0 0 0 0 0 7 0 0 -1 0 0 0 2 21 0 12 0 6 0 0 12 4 0 0 0 11 9 0 0 0 0 -1 0 74 4 10 12 0 0 0 0 10 0 11 0 0 7 0 22 14 0 0 11 13 14 0 11 0 48 5 0 0 18 0 0 0 17 0 9 0 0 0 0 0 0 0 -4 0 10 14 19 0"

# --- TURNS sheet: cumulative totals updated for Week 16 ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 16
$wsTURNS.Range("D3").Value = 12
$wsTURNS.Range("E3").Value = 15

# --- PEN sheet: cumulative totals updated for Week 16 ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 41
